# Rename the worksheet from "Sheet1" to "DoesNotMatter".
# Excel automatically keeps every reference to the sheet name in sync,
# including the <sheets> entry and the _FilterDatabase defined name
# (Sheet1!$B$2:$J$46 -> DoesNotMatter!$B$2:$J$46).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "DoesNotMatter"
